$wb = $excel.ActiveWorkbook

# Worksheets: 1 = Overview, 2 = zh-cn, 3 = ja-jp
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn = $wb.Worksheets.Item(2)
$wsJaJp = $wb.Worksheets.Item(3)

# Update the "Latest Handoff Datetime" / "Latest Handoff Date" values.
# The Overview sheet and the ja-jp sheet share the same underlying value
# (ja-jp's handoff run), the zh-cn sheet has its own handoff run timestamp.
$wsOverview.Range("G2:G7").Value = "2016-07-06 02:21:42"
$wsJaJp.Range("G2:G7").Value = "2016-07-06 02:21:42"
$wsZhCn.Range("G2:G7").Value = "2016-07-06 02:21:33"

# Strip the hash + locale suffix from the generated xlf file names (column F)
# for both locale sheets, so the report shows the plain target file name.
$wsZhCn.Range("F2").Value = "Microsoft.CSharp.RuntimeBinder.Binder.xlf"
$wsZhCn.Range("F3").Value = "Microsoft.CSharp.RuntimeBinder.CSharpArgumentInfo.xlf"
$wsZhCn.Range("F4").Value = "Microsoft.CSharp.RuntimeBinder.CSharpArgumentInfoFlags.xlf"
$wsZhCn.Range("F5").Value = "Microsoft.CSharp.RuntimeBinder.CSharpBinderFlags.xlf"
$wsZhCn.Range("F6").Value = "Microsoft.CSharp.RuntimeBinder.RuntimeBinderException.xlf"
$wsZhCn.Range("F7").Value = "Microsoft.CSharp.RuntimeBinder.RuntimeBinderInternalCompilerException.xlf"

$wsJaJp.Range("F2").Value = "Microsoft.CSharp.RuntimeBinder.Binder.xlf"
$wsJaJp.Range("F3").Value = "Microsoft.CSharp.RuntimeBinder.CSharpArgumentInfo.xlf"
$wsJaJp.Range("F4").Value = "Microsoft.CSharp.RuntimeBinder.CSharpArgumentInfoFlags.xlf"
$wsJaJp.Range("F5").Value = "Microsoft.CSharp.RuntimeBinder.CSharpBinderFlags.xlf"
$wsJaJp.Range("F6").Value = "Microsoft.CSharp.RuntimeBinder.RuntimeBinderException.xlf"
$wsJaJp.Range("F7").Value = "Microsoft.CSharp.RuntimeBinder.RuntimeBinderInternalCompilerException.xlf"
